$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I and J
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Mirror the formatting of the existing header cell (e.g. H1) onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill data rows 2..34: column I is constant 1, column J mirrors column H
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value()
}
